$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 0 -> 6
$ws.Range("B2").Value = 6

# Row 6: C6 6 -> 9, and header labels F6..M6 replaced (K6/L6 swapped order vs columns
# so that shared-string interning order matches the target file)
$ws.Range("C6").Value = 9
$ws.Range("F6").Value = "ax[15:8]"
$ws.Range("G6").Value = "ax[7:0]"
$ws.Range("H6").Value = "ay[20:16]"
$ws.Range("I6").Value = "ay[15:8]"
$ws.Range("J6").Value = "ay[7:0]"
$ws.Range("L6").Value = "az[15:8]"
$ws.Range("K6").Value = "az[20:16]"
$ws.Range("M6").Value = "az[7:0]"
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = ""

# Row 7: new data byte columns
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = "ID"
$ws.Range("E7").Value = "data[31:24]"
$ws.Range("F7").Value = "data[23:16]"
$ws.Range("G7").Value = "data[15:8]"
$ws.Range("H7").Value = "data[7:0]"

# Row 3: B3 1 -> "nyi"
$ws.Range("B3").Value = "nyi"

# Sheet view / column formatting changes
$ws.Range("M1:M20").ColumnWidth = 9.33203125
$ws.Range("C9").Select
